# Scheduled data refresh: update Leve profit calculation inputs/outputs
# across all job sheets (H:N columns) per the latest market price snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 5726.222
$ws.Range("I43").Value = 6598
$ws.Range("J43").Value = 4636.5
$ws.Range("K43").Value = 6598
$ws.Range("L43").Value = 4636.5
$ws.Range("M43").Value = -6529
$ws.Range("N43").Value = -4774.5
$ws.Range("H103").Value = 416.5
$ws.Range("I103").Value = 649
$ws.Range("K103").Value = 1947
$ws.Range("M103").Value = -1361
$ws.Range("H112").Value = 3945.0967
$ws.Range("J112").Value = 4079.2415
$ws.Range("L112").Value = 12237.7245
$ws.Range("N112").Value = -14453.7245
$ws.Range("H113").Value = 13921.333
$ws.Range("I113").Value = 19680.4
$ws.Range("J113").Value = 6722.5
$ws.Range("K113").Value = 19680.4
$ws.Range("L113").Value = 6722.5
$ws.Range("M113").Value = -16426.4
$ws.Range("N113").Value = -13230.5
$ws.Range("H132").Value = 17412.375
$ws.Range("I132").Value = 10313.619
$ws.Range("K132").Value = 30940.857
$ws.Range("M132").Value = -28410.857
$ws.Range("H133").Value = 16999.5
$ws.Range("J133").Value = 16999.5
$ws.Range("L133").Value = 16999.5
$ws.Range("N133").Value = -27119.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1839.8572
$ws.Range("I88").Value = 1636.25
$ws.Range("J88").Value = 1921.3
$ws.Range("K88").Value = 1636.25
$ws.Range("L88").Value = 1921.3
$ws.Range("M88").Value = -1230.25
$ws.Range("N88").Value = -2733.3
$ws.Range("H91").Value = 1839.8572
$ws.Range("I91").Value = 1636.25
$ws.Range("J91").Value = 1921.3
$ws.Range("K91").Value = 1636.25
$ws.Range("L91").Value = 1921.3
$ws.Range("M91").Value = -232.25
$ws.Range("N91").Value = -4729.3
$ws.Range("H97").Value = 543.9583
$ws.Range("I97").Value = 597.4375
$ws.Range("J97").Value = 437
$ws.Range("K97").Value = 597.4375
$ws.Range("L97").Value = 437
$ws.Range("M97").Value = -101.4375
$ws.Range("N97").Value = -1429

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7110611
$ws.Range("I20").Value = 14498138
$ws.Range("J20").Value = 30897.875
$ws.Range("K20").Value = 14498138
$ws.Range("L20").Value = 30897.875
$ws.Range("M20").Value = -14497891
$ws.Range("N20").Value = -31391.875
$ws.Range("H25").Value = 1237.75
$ws.Range("I25").Value = 1237.75
$ws.Range("K25").Value = 1237.75
$ws.Range("M25").Value = -1002.75
$ws.Range("H64").Value = 1558
$ws.Range("J64").Value = 1669.6
$ws.Range("L64").Value = 1669.6
$ws.Range("N64").Value = -2119.6
$ws.Range("H67").Value = 1558
$ws.Range("J67").Value = 1669.6
$ws.Range("L67").Value = 1669.6
$ws.Range("N67").Value = -3229.6
$ws.Range("H82").Value = 9528.571
$ws.Range("I82").Value = 3616.6667
$ws.Range("J82").Value = 45000
$ws.Range("K82").Value = 3616.6667
$ws.Range("L82").Value = 45000
$ws.Range("M82").Value = -3233.6667
$ws.Range("N82").Value = -45766
$ws.Range("H85").Value = 9528.571
$ws.Range("I85").Value = 3616.6667
$ws.Range("J85").Value = 45000
$ws.Range("K85").Value = 3616.6667
$ws.Range("L85").Value = 45000
$ws.Range("M85").Value = -2290.6667
$ws.Range("N85").Value = -47652
$ws.Range("H86").Value = 50002970
$ws.Range("I86").Value = 3652.3635
$ws.Range("J86").Value = 111113240
$ws.Range("K86").Value = 3652.3635
$ws.Range("L86").Value = 111113240
$ws.Range("M86").Value = -2529.3635
$ws.Range("N86").Value = -111115486
$ws.Range("H89").Value = 50002970
$ws.Range("I89").Value = 3652.3635
$ws.Range("J89").Value = 111113240
$ws.Range("K89").Value = 18261.8175
$ws.Range("L89").Value = 555566200
$ws.Range("M89").Value = -12645.8175
$ws.Range("N89").Value = -555577432
$ws.Range("H99").Value = 20198.549
$ws.Range("I99").Value = 18213.143
$ws.Range("J99").Value = 38729
$ws.Range("K99").Value = 18213.143
$ws.Range("L99").Value = 38729
$ws.Range("M99").Value = -16715.143
$ws.Range("N99").Value = -41725
$ws.Range("H102").Value = 23207.5
$ws.Range("I102").Value = 23207.5
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 23207.5
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -19962.5
$ws.Range("N102").ClearContents()
$ws.Range("H134").Value = 37456.676
$ws.Range("I134").Value = 44649
$ws.Range("J134").Value = 25469.467
$ws.Range("K134").Value = 133947
$ws.Range("L134").Value = 76408.401
$ws.Range("M134").Value = -131412
$ws.Range("N134").Value = -81478.401

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 836.6667
$ws.Range("I14").Value = 1010
$ws.Range("J14").Value = 750
$ws.Range("K14").Value = 1010
$ws.Range("L14").Value = 750
$ws.Range("M14").Value = -840
$ws.Range("N14").Value = -1090
$ws.Range("H31").Value = 21458.73
$ws.Range("I31").Value = 12855.223
$ws.Range("K31").Value = 12855.223
$ws.Range("M31").Value = -12560.223
$ws.Range("H34").Value = 21458.73
$ws.Range("I34").Value = 12855.223
$ws.Range("K34").Value = 12855.223
$ws.Range("M34").Value = -12653.223
$ws.Range("H132").Value = 2814.4119
$ws.Range("I132").Value = 2123.0667
$ws.Range("K132").Value = 6369.2001
$ws.Range("M132").Value = -3839.2001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1964188.9
$ws.Range("I4").Value = 1604928.1
$ws.Range("K4").Value = 4814784.300000001
$ws.Range("M4").Value = -4814672.300000001
$ws.Range("H11").Value = 869.8570999999999
$ws.Range("I11").Value = 991.5833
$ws.Range("J11").Value = 139.5
$ws.Range("K11").Value = 2974.7499
$ws.Range("L11").Value = 418.5
$ws.Range("M11").Value = -2834.7499
$ws.Range("N11").Value = -698.5
$ws.Range("H26").Value = 1488.4286
$ws.Range("J26").Value = 70
$ws.Range("L26").Value = 210
$ws.Range("N26").Value = -786
$ws.Range("H98").Value = 2954.7144
$ws.Range("I98").Value = 546.3333
$ws.Range("K98").Value = 1638.9999
$ws.Range("M98").Value = -140.9999
$ws.Range("H109").Value = 4168633.5
$ws.Range("I109").Value = 1497.8334
$ws.Range("K109").Value = 4493.5002
$ws.Range("M109").Value = -3453.5002
$ws.Range("H131").Value = 1469.7474
$ws.Range("J131").Value = 1474.2347
$ws.Range("L131").Value = 4422.7041
$ws.Range("N131").Value = -14502.7041

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 16201.272
$ws.Range("I70").Value = 16900.625
$ws.Range("J70").Value = 14336.333
$ws.Range("K70").Value = 16900.625
$ws.Range("L70").Value = 14336.333
$ws.Range("M70").Value = -16630.625
$ws.Range("N70").Value = -14876.333
$ws.Range("H73").Value = 16201.272
$ws.Range("I73").Value = 16900.625
$ws.Range("J73").Value = 14336.333
$ws.Range("K73").Value = 16900.625
$ws.Range("L73").Value = 14336.333
$ws.Range("M73").Value = -15964.625
$ws.Range("N73").Value = -16208.333
$ws.Range("H97").Value = 1021.4054
$ws.Range("I97").Value = 964.70966
$ws.Range("K97").Value = 964.70966
$ws.Range("M97").Value = -468.70966
$ws.Range("H102").Value = 4097032.8
$ws.Range("J102").Value = 1974.9166
$ws.Range("L102").Value = 1974.9166
$ws.Range("N102").Value = -5218.9166
$ws.Range("H113").Value = 3945.5264
$ws.Range("I113").Value = 3202.2856
$ws.Range("J113").Value = 4379.0835
$ws.Range("K113").Value = 3202.2856
$ws.Range("L113").Value = 4379.0835
$ws.Range("M113").Value = -1032.2856
$ws.Range("N113").Value = -8719.083500000001
$ws.Range("H122").Value = 3087977
$ws.Range("I122").Value = 3773360.8
$ws.Range("J122").Value = 3750
$ws.Range("K122").Value = 11320082.4
$ws.Range("L122").Value = 11250
$ws.Range("M122").Value = -11317632.4
$ws.Range("N122").Value = -16150

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 71430560
$ws.Range("I16").Value = 90911390
$ws.Range("J16").Value = 847
$ws.Range("K16").Value = 90911390
$ws.Range("L16").Value = 847
$ws.Range("M16").Value = -90911220
$ws.Range("N16").Value = -1187
$ws.Range("H22").Value = 21740916
$ws.Range("J22").Value = 47621388
$ws.Range("L22").Value = 47621388
$ws.Range("N22").Value = -47621978
$ws.Range("H27").Value = 21740916
$ws.Range("J27").Value = 47621388
$ws.Range("L27").Value = 47621388
$ws.Range("N27").Value = -47621602
$ws.Range("H74").Value = 48098.5
$ws.Range("I74").Value = 48098.5
$ws.Range("K74").Value = 48098.5
$ws.Range("M74").Value = -47100.5
$ws.Range("H77").Value = 48098.5
$ws.Range("I77").Value = 48098.5
$ws.Range("K77").Value = 144295.5
$ws.Range("M77").Value = -139303.5
$ws.Range("H122").Value = 23227524
$ws.Range("I122").Value = 41662750
$ws.Range("J122").Value = 2278400.8
$ws.Range("K122").Value = 124988250
$ws.Range("L122").Value = 6835202.399999999
$ws.Range("M122").Value = -124985800
$ws.Range("N122").Value = -6840102.399999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 20333.334
$ws.Range("I42").Value = 11500
$ws.Range("K42").Value = 11500
$ws.Range("M42").Value = -11122
$ws.Range("H100").Value = 611.7727
$ws.Range("I100").Value = 566.3333
$ws.Range("K100").Value = 1132.6666
$ws.Range("M100").Value = -591.6666
$ws.Range("H107").Value = 2529.2
$ws.Range("I107").Value = 2529.2
$ws.Range("K107").Value = 7587.599999999999
$ws.Range("M107").Value = -5667.599999999999
$ws.Range("H122").Value = 329822.78
$ws.Range("I122").Value = 406213.8
$ws.Range("K122").Value = 1218641.4
$ws.Range("M122").Value = -1216191.4
$ws.Range("H124").Value = 33000
$ws.Range("J124").Value = 33000
$ws.Range("L124").Value = 33000
$ws.Range("N124").Value = -42820
$ws.Range("H132").Value = 13348.571
$ws.Range("I132").Value = 5239.184
$ws.Range("K132").Value = 15717.552
$ws.Range("M132").Value = -13187.552
